# "added batch fft and classification to app"
#
# The "frequency" sheet (sheet1) holds raw FFT-trial measurements in C4:J5
# (row 4 = window 0.15, row 5 = window 0.2). Re-running the batch FFT /
# classification step produced new trial numbers for most of those cells;
# C8:J8 (AVERAGE / STDEV.S over C4:D5, E4:F5, G4:H5, I4:J5) are formulas and
# recalculate automatically once the inputs change. The sheet's selection
# also moved from M6 to K5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("frequency")

# Row 4 (window = 0.15)
$ws.Range("C4").Value = 1.9851180881350801
$ws.Range("D4").Value = 1.9789035803178201
# E4 (1.54979062551807) is unchanged
$ws.Range("F4").Value = 1.478073722682
$ws.Range("G4").Value = 1.0041127468914399
$ws.Range("H4").Value = 0.96838985694941304
$ws.Range("I4").Value = 0.41040337021083301
$ws.Range("J4").Value = 0.68272611380339998

# Row 5 (window = 0.2)
$ws.Range("C5").Value = 2.0269129436259301
# D5 (2.00587260674684) is unchanged
$ws.Range("E5").Value = 1.5262150524318701
$ws.Range("F5").Value = 1.5054790951119701
$ws.Range("G5").Value = 1.00535569580512
$ws.Range("H5").Value = 0.97550431478145905
$ws.Range("I5").Value = 1.0081890395142601
$ws.Range("J5").Value = 0.50540244401992596

# Selection moved from M6 to K5
$ws.Activate() | Out-Null
$ws.Range("K5").Select() | Out-Null
